# Update cryptos list: refreshed Price (D) / Volume(1h) (E) snapshot values,
# plus two coins (rows 33/34 and 39/40) swapping rank order.
#
# NOTE: Price strings (column D) are stored as *text* in the workbook (the
# values use '.' as a thousands separator, e.g. "27.950.00", so they are not
# valid numbers). Excel auto-converts a plain numeric-looking string typed
# into .Value back into a float, which would both corrupt values like
# "212.73" (-> 212.73 the number, dropping formatting) and silently retype
# the cell. Prefixing with a leading apostrophe forces Excel to keep it as
# literal text, matching the original text/shared-string cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.950.00"
$ws.Range("E2").Value = "  -0.03%  "

$ws.Range("D3").Value = "'1.639.88"
$ws.Range("E3").Value = "  -0.27%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").Value = "'212.73"
$ws.Range("E5").Value = "  +0.16%  "

$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").Value = "'23.44"
$ws.Range("E8").Value = "  -0.50%  "

$ws.Range("E9").Value = "  -2.24%  "

$ws.Range("D10").Value = "'0.0614"
$ws.Range("E10").Value = "  +0.14%  "

$ws.Range("D11").Value = "'0.0883"
$ws.Range("E11").Value = "  +1.81%  "

$ws.Range("D12").Value = "'1.873.75"
$ws.Range("E12").Value = "  -0.13%  "

$ws.Range("D13").Value = "'1.642.49"
$ws.Range("E13").Value = "  -0.09%  "

$ws.Range("D14").Value = "'4.08"
$ws.Range("E14").Value = "  +0.42%  "

$ws.Range("E15").Value = "  +0.94%  "

$ws.Range("D16").Value = "'65.48"
$ws.Range("E16").Value = "  -0.25%  "

$ws.Range("D17").Value = "'27.948.53"
$ws.Range("E17").Value = "  +0.06%  "

$ws.Range("D18").Value = "'232.70"
$ws.Range("E18").Value = "  +0.10%  "

$ws.Range("D19").Value = "'0.0₃0722"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").Value = "'7.58"
$ws.Range("E20").Value = "  -1.23%  "

$ws.Range("E21").Value = "  +0.22%  "

$ws.Range("D22").Value = "'10.53"
$ws.Range("E22").Value = "  -2.08%  "

$ws.Range("E23").Value = "  -0.64%  "

$ws.Range("E24").Value = "  -4.14%  "

$ws.Range("D25").Value = "'153.07"
$ws.Range("E25").Value = "  +1.52%  "

$ws.Range("D26").Value = "'6.89"
$ws.Range("E26").Value = "  -0.56%  "

$ws.Range("D27").Value = "'15.66"
$ws.Range("E27").Value = "  -0.28%  "

$ws.Range("E28").Value = "  -0.42%  "

$ws.Range("E29").Value = "  +0.22%  "

$ws.Range("E30").Value = "  +0.24%  "

$ws.Range("D31").Value = "'0.0484"
$ws.Range("E31").Value = "  +0.32%  "

$ws.Range("E32").Value = "  +2.69%  "

# Rows 33/34 swap: InternetComputer(DFINITY) now ranks above Maker.
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'3.09"
$ws.Range("E33").Value = "  -0.30%  "

$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "'1.405.58"
$ws.Range("E34").Value = "  -4.28%  "

$ws.Range("E35").Value = "  +1.05%  "

$ws.Range("E36").Value = "  +1.71%  "

$ws.Range("E37").Value = "  +0.59%  "

$ws.Range("E38").Value = "  +0.55%  "

# Rows 39/40 swap: ARBITRUM now ranks above TrustWalletToken.
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'0.879"
$ws.Range("E39").Value = "  -0.94%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'0.926"
$ws.Range("E40").Value = "  -0.65%  "

$ws.Range("E41").Value = "  +0.64%  "

$ws.Range("E42").Value = "  +0.16%  "

$ws.Range("E43").Value = "  +4.86%  "

$ws.Range("D44").Value = "'67.04"
$ws.Range("E44").Value = "  -3.18%  "

$ws.Range("E45").Value = "  +2.46%  "

$ws.Range("E46").Value = "  -0.44%  "

$ws.Range("D47").Value = "'1.781.58"
$ws.Range("E47").Value = "  -0.58%  "

$ws.Range("D48").Value = "'87.87"
$ws.Range("E48").Value = "  +0.00%  "

$ws.Range("E49").Value = "  +2.15%  "

$ws.Range("D50").Value = "'0.1000"
$ws.Range("E50").Value = "  -0.36%  "

$ws.Range("E51").Value = "  +0.08%  "

Write-Host "Updated cryptos worksheet"
